$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - sheet1
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 298
$ws.Range("F4").Value = 40
$ws.Range("F5").Value = 47
$ws.Range("F7").Value = 1244
$ws.Range("F8").Value = 386
$ws.Range("F9").Value = 244
$ws.Range("F10").Value = 357
$ws.Range("F11").Value = 8244
$ws.Range("F13").Value = 9967
$ws.Range("F14").Value = 84
$ws.Range("F18").Value = 495
$ws.Range("F27").Value = 1707
$ws.Range("F28").Value = 45
$ws.Range("F30").Value = 317
$ws.Range("F31").Value = 275
$ws.Range("F33").Value = 552
$ws.Range("F34").Value = 1022
$ws.Range("F37").Value = 1403
$ws.Range("F38").Value = 420
$ws.Range("F42").Value = 488
$ws.Range("F43").Value = 314
$ws.Range("F46").Value = 109
$ws.Range("F48").Value = 31
$ws.Range("F49").Value = 36

# Sheet "演出" (Performance) - sheet2
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F15").Value = 56
$ws.Range("F19").Value = 31

# Sheet "本地生活" (Local life) - sheet3
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2768
$ws.Range("F5").Value = 195

# Sheet "全部类型" (All types) - sheet4
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 298
$ws.Range("F6").Value = 195
$ws.Range("F7").Value = 40
$ws.Range("F8").Value = 47
$ws.Range("F10").Value = 1244
$ws.Range("F11").Value = 386
$ws.Range("F14").Value = 244
$ws.Range("F16").Value = 8244
$ws.Range("F18").Value = 9968
$ws.Range("F24").Value = 1707
$ws.Range("F25").Value = 45
$ws.Range("F26").Value = 317
$ws.Range("F27").Value = 275
$ws.Range("F30").Value = 552
$ws.Range("F36").Value = 1403
$ws.Range("F37").Value = 420
$ws.Range("F38").Value = 56
$ws.Range("F41").Value = 488
$ws.Range("F42").Value = 314
$ws.Range("F45").Value = 31
$ws.Range("F48").Value = 31
$ws.Range("F49").Value = 36
